$d = $word.ActiveDocument

# Update the due date text: "November 9th, 2017" -> "November 8th, 2018"
$d.Content.Find.Execute("November 9", $true, $false, $false, $false, $false,
                         $true, 1, $false, "November 8", 2)
$d.Content.Find.Execute("2017", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2018", 2)

# Word keeps a hidden "_GoBack" bookmark marking the location of the most recent
# edit. Before this change it sat at the end of the first ("Assignment / COS318")
# paragraph; since the edit happened on the due-date line, move it there - right
# after the newly-typed year, before the line break that precedes
# "Turn in all files...".
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$dueDateParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Due Date*") {
        $dueDateParagraph = $p
        break
    }
}

$editRange = $dueDateParagraph.Range
$editRange.Find.Execute(", 2018", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
$editRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $editRange)
